$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (2023-09-02 -> 2023-09-03) for every data row (rows 2-223).
$ws.Range("C2:C223").Value = 45172
